# Applies the "Correcciones general y test" edit:
#  - Inserts a new row "user | role_user | required | 1" right after the
#    existing "user | last_session | display_item | none" row (old row 16).
#  - Inserts a new row "menu | menu_type | display_list | show" right after
#    the existing "menu | parent_id | display_list | none" row (which, after
#    the first insertion, has shifted from row 19 to row 20).
#  - Leaves every other row untouched; Excel itself takes care of shifting
#    row numbers / the <dimension> ref when rows are inserted.
#  - Restores the final selection to D18, matching the authored diff.
#
# New shared-string cells must be written in the same order they get a new
# unique-string slot in the saved workbook: "menu_type" has to land on index
# 72 and "role_user" on index 73, so we type the menu_type cell before the
# role_user cell even though role_user's row is physically earlier.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("edit-fields")

# --- Insert "menu | menu_type | display_list | show" ----------------------
# Today this lands right after "menu | parent_id | display_list | none"
# (row 19) and before "menu | order | display_list | show" (row 20).
$ws.Rows.Item(20).Insert()
$ws.Range("B20").Value = "menu_type"
$ws.Range("A20").Value = "menu"
$ws.Range("C20").Value = "display_list"
$ws.Range("D20").Value = "show"

# --- Insert "user | role_user | required | 1" ------------------------------
# This lands right after "user | last_session | display_item | none"
# (row 16) and before "menu | level | preset | 1" (row 17).
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "user"
$ws.Range("B17").Value = "role_user"
$ws.Range("C17").Value = "required"
$ws.Range("D17").Value = 1

# --- Restore the reported selection ---------------------------------------
$ws.Range("D18").Select()
